# Generate Report for Handoff
#
# Semantics of this edit (per the target diff):
#  - The handoff file fdf07564-fe57-4b5e-8e0c-798793e39fcd.md was replaced by a
#    newly generated handoff package 447830c7-cebb-489f-ba28-ec609559dc98.md
#    (and its paired .xlf deliverables), with refreshed handoff timestamps.
#  - A brand-new source file ffff865736d8-fcea-42ea-a321-97636be0831d.md showed
#    up in the same handoff and is now tracked as an additional row on every
#    sheet (Overview / zh-cn / de-de).
#
# We rebuild hyperlinks from scratch on every sheet (Hyperlinks.Delete() then
# Hyperlinks.Add(...)) instead of mutating existing Hyperlink objects in place,
# because mutating an existing Hyperlinks.Item(...) in this host duplicates the
# <hyperlink> entry rather than replacing it.

$wb = $excel.ActiveWorkbook

$oldGuid = "fdf07564-fe57-4b5e-8e0c-798793e39fcd"
$newGuid = "447830c7-cebb-489f-ba28-ec609559dc98"
$newGuid2 = "ffff865736d8-fcea-42ea-a321-97636be0831d"
$oldHash = "4a5a64796d3e740df5ee133c6ce91e66321171ee"
$newHash = "f28090322b799169d72c656f0eb2c2578504b8d6"

$newMd = "$newGuid.md"
$newMd2 = "$newGuid2.md"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$handoffDate = "2016-03-23 04:43:06"
$zhHandoffDatetime = "2016-03-23 04:42:56"
$deHandoffDatetime = "2016-03-23 04:43:06"
$epoch = "0001-01-01 00:00:00"
$readyStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $handoffDate

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = $readyStatus
$ws.Range("C3").Value = $readyStatus
$ws.Range("D3").Value = $handoffDate
$ws.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd2", "", "", $newMd2) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newZhXlf
$ws.Range("E2").Value = $zhHandoffDatetime

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $readyStatus
$ws.Range("D3").Value = $newZhXlf
$ws.Range("E3").Value = $zhHandoffDatetime
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = $epoch
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca1f452ffe66e3917bac9a84b219a215c4effb8e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$newZhXlf", "", "", $newZhXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd2", "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ca1f452ffe66e3917bac9a84b219a215c4effb8e/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$newZhXlf", "", "", $newZhXlf) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()

$ws.Range("A2").Value = $newMd
$ws.Range("D2").Value = $newDeXlf
$ws.Range("E2").Value = $deHandoffDatetime

$ws.Range("A3").Value = $newMd2
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $readyStatus
$ws.Range("D3").Value = $newDeXlf
$ws.Range("E3").Value = $deHandoffDatetime
$ws.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H3").Value = $epoch
$ws.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("J3").Value = "Include"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd", "", "", $newMd) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c4d0396c02989de096852f97d720b992ccfefb2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$newDeXlf", "", "", $newDeXlf) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5884ce024b9e3ce76ee84248c08110dcdd726655/e2e/$newMd2", "", "", $newMd2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9c4d0396c02989de096852f97d720b992ccfefb2/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$newDeXlf", "", "", $newDeXlf) | Out-Null

Write-Host "Applied handoff report regeneration edits."
